# excel_mac_2011-formatting.xlsx edit
#
# - Rename Sheet1 -> "first sheet"
# - Add a new "Sheet2" after it, with A1 = "second sheet"
# - Sheet1: add a far-away header cell (AA1 = "this far column")
# - Sheet1: restyle the B1/C1 header cells (big colored font w/ border,
#   and a bordered+filled+centered/wrapped header cell)
# - Update selections on both sheets, keep Sheet1 as the active tab

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: rename ---------------------------------------------------
$ws1.Name = "first sheet"

# --- Sheet1: new far cell (written first so its shared-string index
#     comes out before "second sheet", matching the original edit) ------
$ws1.Range("AA1").Value = "this far column"

# --- Add Sheet2 right after Sheet1 -------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws2.Range("A1").Value = "second sheet"

# --- Sheet1: B1 header cell -> big colored font + full thin border -----
$b1 = $ws1.Range("B1")
$b1.Font.Size = 28
$b1.Font.ColorIndex = 26
$b1.Borders.LineStyle = 1
$b1.Borders.Weight = 2

# --- Sheet1: C1 header cell -> fill + border + centered/wrapped text ---
$c1 = $ws1.Range("C1")
$c1.Font.ColorIndex = 1
$c1.Interior.ColorIndex = 43
$c1.Interior.Pattern = 1
$c1.Borders(7).LineStyle = 1
$c1.Borders(7).Weight = 2
$c1.HorizontalAlignment = -4108
$c1.VerticalAlignment = -4108
$c1.WrapText = $true

# --- Sheet2: selection (doesn't move the active tab) -------------------
$ws2.Range("C5").Select()

# --- Sheet1: stays the active tab, with C1 selected ---------------------
$ws1.Activate()
$ws1.Range("C1").Select()
